# upgrade YAML file from beta to V1.
#
# The only content-level change recorded by the diff is that the hidden
# "_GoBack" bookmark (Word's "last edit position" marker) moved from the
# empty, shaded paragraph right after the "Dashboard" heading to the
# empty paragraph right before it. Re-adding a bookmark with a name that
# already exists elsewhere in the document moves it (Word enforces unique
# bookmark names), so a single Bookmarks.Add call reproduces both the
# removal of the old <w:bookmarkStart/><w:bookmarkEnd/> pair and the
# insertion of the new one.

$d = $word.ActiveDocument
$paras = $d.Paragraphs

$targetIndex = -1
$i = 0
foreach ($p in $paras) {
    $i = $i + 1
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "Dashboard" -and $p.Style.NameLocal -eq "Heading 2") {
        $targetIndex = $i
    }
}

if ($targetIndex -gt 1) {
    $prevPara = $paras.Item($targetIndex - 1)
    $d.Bookmarks.Add("_GoBack", $prevPara.Range)
    Write-Output "Moved _GoBack bookmark to paragraph $($targetIndex - 1) (before 'Dashboard')."
} else {
    Write-Output "Could not locate 'Dashboard' heading; no changes made."
}
